$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 97, shifting existing rows 97-100 down to 99-102
$ws.Rows("97:98").Insert()

# Fill new row 97 (Española, 30 unidades)
$ws.Cells.Item(97, 1).Value = 10
$ws.Cells.Item(97, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(97, 3).Value = "La Araucanía"
$ws.Cells.Item(97, 4).Value = 44448
$ws.Cells.Item(97, 5).Value = 9
$ws.Cells.Item(97, 6).Value = 100112013
$ws.Cells.Item(97, 7).Value = "Alcachofa"
$ws.Cells.Item(97, 8).Value = "Española"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 65
$ws.Cells.Item(97, 11).Value = 14000
$ws.Cells.Item(97, 12).Value = 14000
$ws.Cells.Item(97, 13).Value = 14000
$ws.Cells.Item(97, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(97, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(97, 16).Value = 467
$ws.Cells.Item(97, 17).Value = 30
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# Fill new row 98 (Madrigal, 40 unidades)
$ws.Cells.Item(98, 1).Value = 10
$ws.Cells.Item(98, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(98, 3).Value = "La Araucanía"
$ws.Cells.Item(98, 4).Value = 44448
$ws.Cells.Item(98, 5).Value = 9
$ws.Cells.Item(98, 6).Value = 100112013
$ws.Cells.Item(98, 7).Value = "Alcachofa"
$ws.Cells.Item(98, 8).Value = "Madrigal"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 65
$ws.Cells.Item(98, 11).Value = 14000
$ws.Cells.Item(98, 12).Value = 14000
$ws.Cells.Item(98, 13).Value = 14000
$ws.Cells.Item(98, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(98, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(98, 16).Value = 350
$ws.Cells.Item(98, 17).Value = 40
$ws.Cells.Item(98, 18).Value = "Hortaliza"

Write-Host "done"
